$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "66.444.04"
$ws.Range("E2").Value = "  -3.52%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.560.45"
$ws.Range("E3").Value = "  -3.84%  "
$ws.Range("E4").Value = "  +0.11%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "579.16"
$ws.Range("E5").Value = "  -5.42%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "186.13"
$ws.Range("E6").Value = "  -0.57%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "3.555.85"
$ws.Range("E7").Value = "  -3.93%  "
$ws.Range("E8").Value = "  -3.20%  "
$ws.Range("E9").Value = "  +0.14%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.666"
$ws.Range("E10").Value = "  -6.10%  "
$ws.Range("E11").Value = "  -8.88%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "52.71"
$ws.Range("E12").Value = "  -4.94%  "
$ws.Range("E13").Value = "  -9.97%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "9.75"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "4.125.89"
$ws.Range("E15").Value = "  -3.87%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.559.64"
$ws.Range("E16").Value = "  -3.90%  "
$ws.Range("E17").Value = "  -0.87%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "18.22"
$ws.Range("E18").Value = "  -4.93%  "
$ws.Range("B19").Value = "WrappedBTC"
$ws.Range("C19").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "66.292.23"
$ws.Range("E19").Value = "  -3.41%  "
$ws.Range("B20").Value = "Uniswap"
$ws.Range("C20").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "12.11"
$ws.Range("E20").Value = "  -5.75%  "
$ws.Range("E21").Value = "  -6.17%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "392.90"
$ws.Range("E22").Value = "  -3.83%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.30"
$ws.Range("E23").Value = "  -5.64%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "85.49"
$ws.Range("E24").Value = "  -3.61%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "11.07"
$ws.Range("E25").Value = "  +2.06%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.88"
$ws.Range("E26").Value = "  -4.30%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "12.35"
$ws.Range("E27").Value = "  -2.57%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "6.04"
$ws.Range("E28").Value = "  -0.01%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "3.51"
$ws.Range("E29").Value = "  -5.98%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "8.86"
$ws.Range("E30").Value = "  -7.44%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "30.93"
$ws.Range("E31").Value = "  -5.84%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "7.06"
$ws.Range("E32").Value = "  -1.32%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "625.19"
$ws.Range("E33").Value = "  +0.35%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "12.13"
$ws.Range("E34").Value = "  -3.18%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "63.62"
$ws.Range("E35").Value = "  -2.92%  "
$ws.Range("E36").Value = "  -6.67%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "41.11"
$ws.Range("E37").Value = "  -6.21%  "
$ws.Range("E38").Value = "  +0.30%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.395"
$ws.Range("E39").Value = "  -3.15%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0₃0765"
$ws.Range("E40").Value = "  -4.29%  "
$ws.Range("E41").Value = "  -5.93%  "
$ws.Range("E42").Value = "  -0.17%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.966.94"
$ws.Range("E43").Value = "  +3.80%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.80"
$ws.Range("E44").Value = "  -6.50%  "
$ws.Range("E45").Value = "  -3.96%  "
$ws.Range("E46").Value = "  -7.46%  "
$ws.Range("E47").Value = "  -6.47%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "3.08"
$ws.Range("E48").Value = "  -0.05%  "
$ws.Range("B49").Value = "Monero"
$ws.Range("C49").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "138.00"
$ws.Range("E49").Value = "  -2.45%  "
$ws.Range("B50").Value = "THORChain"
$ws.Range("C50").Value = "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "8.43"
$ws.Range("E50").Value = "  -6.45%  "
$ws.Range("B51").Value = "Stacks"
$ws.Range("C51").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.72"
$ws.Range("E51").Value = "  -0.90%  "
